$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 817
$ws.Range("F7").Value = 393
$ws.Range("F8").Value = 4649
$ws.Range("F9").Value = 4649
$ws.Range("F11").Value = 118
$ws.Range("F12").Value = 152
$ws.Range("F15").Value = 107
$ws.Range("F16").Value = 7351
$ws.Range("F17").Value = 248
$ws.Range("F21").Value = 506
$ws.Range("F22").Value = 1337
$ws.Range("F23").Value = 439
$ws.Range("F24").Value = 6278
$ws.Range("E25").Value = "2024.04.20 10:00-04.21 17:00"
$ws.Range("F28").Value = 6139
$ws.Range("F29").Value = 139
$ws.Range("F34").Value = 6356
$ws.Range("F35").Value = 21
$ws.Range("F36").Value = 201
$ws.Range("F37").Value = 95
$ws.Range("F38").Value = 21
$ws.Range("F43").Value = 53
$ws.Range("F46").Value = 419
$ws.Range("F47").Value = 2127
$ws.Range("F48").Value = 39
$ws.Range("F49").Value = 1069

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 126
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = 45
$ws.Range("F6").Value = 118
$ws.Range("F8").Value = 9

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 45
$ws.Range("F8").Value = 393
$ws.Range("F9").Value = 4649
$ws.Range("F10").Value = 4649
$ws.Range("F12").Value = 118
$ws.Range("F13").Value = 152
$ws.Range("F16").Value = 107
$ws.Range("F17").Value = 7351
$ws.Range("F18").Value = 248
$ws.Range("F20").Value = 506
$ws.Range("F21").Value = 1337
$ws.Range("F22").Value = 118
$ws.Range("F23").Value = 6278
$ws.Range("E24").Value = "2024.04.20 10:00-04.21 17:00"
$ws.Range("F27").Value = 9
$ws.Range("F29").Value = 6139
$ws.Range("F30").Value = 139
$ws.Range("F36").Value = 6356
$ws.Range("F37").Value = 21
$ws.Range("F38").Value = 201
$ws.Range("F39").Value = 95
$ws.Range("F46").Value = 419
$ws.Range("F48").Value = 2127
$ws.Range("F49").Value = 39
